$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $escaped = $text.Replace('"', '""')
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "66.551.43"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.513.47"
$ws.Range("E3").Value = "  -4.84%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "583.52"
$ws.Range("E5").Value = "  -2.05%  "
Set-TextValue "D6" "171.26"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.81%  "
$ws.Range("D9").Value = "2.512.72"
$ws.Range("E9").Value = "  -4.83%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -4.62%  "
$ws.Range("E13").Value = "  -2.29%  "
Set-TextValue "D14" "26.72"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "2.971.15"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("D17").Value = "66.301.49"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "2.509.94"
$ws.Range("E18").Value = "  -3.93%  "
Set-TextValue "D19" "7.81"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("E20").Value = "  -6.74%  "
Set-TextValue "D21" "347.90"
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("E22").Value = "  -3.28%  "
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  -0.07%  "
Set-TextValue "D26" "70.00"
$ws.Range("E26").Value = "  -0.11%  "
Set-TextValue "D27" "9.95"
$ws.Range("E27").Value = "  -4.32%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "2.631.52"
$ws.Range("E29").Value = "  -5.14%  "
$ws.Range("D30").Value = "0.0₃0975"
$ws.Range("E30").Value = "  -3.83%  "
Set-TextValue "D31" "523.75"
$ws.Range("E31").Value = "  -5.02%  "
Set-TextValue "D32" "8.09"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("E33").Value = "  -3.37%  "
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -3.28%  "
Set-TextValue "D38" "156.87"
$ws.Range("E38").Value = "  -0.67%  "
Set-TextValue "D39" "18.65"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("E41").Value = "  -3.33%  "
Set-TextValue "D42" "1.80"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +1.71%  "
Set-TextValue "D46" "39.44"
$ws.Range("E46").Value = "  -1.62%  "
Set-TextValue "D47" "149.04"
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("E48").Value = "  -4.41%  "
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "0.0₆0269"
$ws.Range("E51").Value = "  -11.43%  "
